# The XML-mapped table "Tabel2" had two columns whose header names used the
# old, verbose convention (PlatformLogoPath / ImagePath). The commit renames
# them to the shorter names Platform / Image. Renaming the header cells in
# an XML-mapped table cascades through the shared-string table, the sheet
# cell values, and the table column definitions (xl/tables/table1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 = "PlatformLogoPath" -> "Platform"
$ws.Range("D1").Value = "Platform"
# E1 = "ImagePath" -> "Image"
$ws.Range("E1").Value = "Image"

# Move the active selection to H12 (as recorded in the saved view state).
$ws.Range("H12").Select() | Out-Null
